# Generate Report for Handback
#
# Both localized-file rows (a.md, b.md) have come back "in sync" with the
# source, so the per-language sheets move from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Handback DateTime is
# refreshed to the moment of this report run, and the now-stale
# Latest Handback Name / Error Detail fields are cleared out.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------------
# Columns E (zh-cn) / F (de-de) mirror each language sheet's Status column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet ---------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhHandbackTime = "2017-02-22 07:03:26"

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("L2").Value = $zhHandbackTime
$wsZh.Range("L3").Value = $zhHandbackTime

$wsZh.Range("M2").Value = ""
$wsZh.Range("M3").Value = ""

$wsZh.Range("R2").Value = ""
$wsZh.Range("R3").Value = ""

# --- de-de sheet -----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deHandbackTime = "2017-02-22 07:03:49"

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("L2").Value = $deHandbackTime
$wsDe.Range("L3").Value = $deHandbackTime

$wsDe.Range("M2").Value = ""
$wsDe.Range("M3").Value = ""

$wsDe.Range("R2").Value = ""
$wsDe.Range("R3").Value = ""

# --- Column widths: the text-content changes above (shorter status text
# removed, longer status text added, long error-detail text cleared, long
# handback-name text cleared) reflow the autosized columns on all three
# sheets. Re-apply AutoFit to the affected columns so widths track content.
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

foreach ($ws in @($wsZh, $wsDe)) {
    $ws.Columns.Item(3).AutoFit()
    $ws.Columns.Item(13).AutoFit()
    $ws.Columns.Item(18).AutoFit()
}
